# "change file name to english, add summary of collection & iterator"
#
# Rewrites D2 from an integer percent-looking number (50) into a real
# percentage value (0.5, formatted as 0%), and appends rows 3-8 describing
# python "collection & iterator" topics, each with a completion percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (Excel VBA-style BGR-packed decimal)
$green = 5287936   # RGB(0,176,80) -> R + G*256 + B*65536

# --- Row 2: fix D2 to be a real percentage (0.5 formatted as 0%) ---
$ws.Range("D2").Value = 0.5
$ws.Range("D2").NumberFormat = "0%"

# --- Row 3: python / 基础 / 运算符 / 80% ---
$ws.Range("A3").Value = "python"
$ws.Range("A3").Font.Color = $green

$ws.Range("B3").Value = "基础"
$ws.Range("B3").Font.Color = $green
$ws.Range("B3").Font.Family = 3

$ws.Range("C3").Value = "运算符"
$ws.Range("C3").Font.Color = $green
$ws.Range("C3").Font.Family = 3

$ws.Range("D3").Value = 0.8
$ws.Range("D3").Font.Color = $green
$ws.Range("D3").Font.Family = 3
$ws.Range("D3").NumberFormat = "0%"

# --- Row 4: python / 数据结构 / set / 80% ---
$ws.Range("A4").Value = "python"
$ws.Range("A4").Font.Color = $green

$ws.Range("B4").Value = "数据结构"
$ws.Range("B4").Font.Color = $green
$ws.Range("B4").Font.Family = 3

$ws.Range("C4").Value = "set"
$ws.Range("C4").Font.Color = $green
$ws.Range("C4").Font.Family = 3

$ws.Range("D4").Value = 0.8
$ws.Range("D4").Font.Color = $green
$ws.Range("D4").Font.Family = 3
$ws.Range("D4").NumberFormat = "0%"

# --- Row 5: python / 数据结构 / deque / 80% ---
$ws.Range("A5").Value = "python"
$ws.Range("A5").Font.Color = $green

$ws.Range("B5").Value = "数据结构"
$ws.Range("B5").Font.Color = $green
$ws.Range("B5").Font.Family = 3

$ws.Range("C5").Value = "deque"
$ws.Range("C5").Font.Color = $green
$ws.Range("C5").Font.Family = 3

$ws.Range("D5").Value = 0.8
$ws.Range("D5").Font.Color = $green
$ws.Range("D5").Font.Family = 3
$ws.Range("D5").NumberFormat = "0%"

# --- Row 6: python / 基础 / 迭代器 / 80% ---
$ws.Range("A6").Value = "python"
$ws.Range("A6").Font.Color = $green

$ws.Range("B6").Value = "基础"
$ws.Range("B6").Font.Color = $green
$ws.Range("B6").Font.Family = 3

$ws.Range("C6").Value = "迭代器"
$ws.Range("C6").Font.Color = $green
$ws.Range("C6").Font.Family = 3

$ws.Range("D6").Value = 0.8
$ws.Range("D6").Font.Color = $green
$ws.Range("D6").Font.Family = 3
$ws.Range("D6").NumberFormat = "0%"

# --- Row 7: python (green CJK font, matches original workbook) / 基础 / 生成器 / 80% ---
$ws.Range("A7").Value = "python"
$ws.Range("A7").Font.Color = $green
$ws.Range("A7").Font.Family = 3

$ws.Range("B7").Value = "基础"
$ws.Range("B7").Font.Color = $green
$ws.Range("B7").Font.Family = 3

$ws.Range("C7").Value = "生成器"
$ws.Range("C7").Font.Color = $green
$ws.Range("C7").Font.Family = 3

$ws.Range("D7").Value = 0.8
$ws.Range("D7").Font.Color = $green
$ws.Range("D7").Font.Family = 3
$ws.Range("D7").NumberFormat = "0%"

# --- Row 8: python / 数据结构 / ChainMap (default color) / 70% ---
$ws.Range("A8").Value = "python"
$ws.Range("B8").Value = "数据结构"

$ws.Range("C8").Value = "ChainMap"
$ws.Range("C8").Font.Family = 3

$ws.Range("D8").Value = 0.7
$ws.Range("D8").NumberFormat = "0%"

# --- Column C width (closest reachable approximation of the authored
#     best-fit width; this runtime's ColumnWidth setter quantizes to
#     1/7-character steps so the exact 10.33203125 best-fit value from
#     real Excel can't be reproduced bit-for-bit) ---
$ws.Columns("C").ColumnWidth = 9.5

# --- Final selection, matching the authored workbook state ---
[void]$ws.Range("D4").Select()
